# The trailing empty paragraph (<w:p/>) must become a paragraph that
# contains a single, empty run (<w:p><w:r><w:t/></w:r></w:p>).
$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertBefore("")
